$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G: difference between LoudnessThreshold (E) and DiffThreshold (F).
# G2 is entered on its own (plain formula); G3:G13 are filled as one
# operation so they share a single formula definition (shared formula group).
$ws.Range("G2").Formula = "=E2-F2"
$ws.Range("G3:G13").Formula = "=E3-F3"

# Give the new column a sensible width, matching the "best fit" sizing used
# by the other data columns (A, C:F).
$ws.Columns.Item(7).ColumnWidth = 10.1666666666667

# The block of leftover/empty formatted cells below & right of the table
# (rows 14-21, columns H:Q) is trimmed back - only the still-meaningful
# formatting remains, the rest is cleared out entirely.
$ws.Range("N14:Q14").Clear()
$ws.Range("Q15").Clear()
$ws.Range("J16:Q16").Clear()
$ws.Range("P17:Q17").Clear()
$ws.Range("P18:Q18").Clear()
$ws.Range("J19:Q20").Clear()
$ws.Range("J21:K21").Clear()

# Reflect the active selection left after filling in column G.
$ws.Range("G2:G13").Select()
